$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.076.18"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "3.354.97"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Formula = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Formula = "'569.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Formula = "'135.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.353.05"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Formula = "'7.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Formula = "'0.123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "3.922.23"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Formula = "'25.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "3.346.27"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "61.173.27"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Formula = "'13.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Formula = "'5.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Formula = "'376.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("D23").Formula = "'0.551"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").Value = "3.485.64"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Formula = "'70.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Formula = "'0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").Formula = "'1.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.26%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").Formula = "'8.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("D33").Formula = "'2.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D35").Formula = "'23.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Formula = "'5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.42%  "
$ws.Range("D37").Formula = "'6.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Formula = "'164.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("D40").Formula = "'0.0756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Formula = "'1.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Formula = "'0.766"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Formula = "'41.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").Formula = "'4.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("E47").Value = "  -4.88%  "
$ws.Range("D48").Formula = "'23.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").Value = "2.346.61"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -2.04%  "
